$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 27 (TestScenario_26 / TestScenario_26.TestCase_1):
#   Approved/Rejected was "Approved" with no reject reason.
#   It is now re-reviewed and marked "Rejected" with reason "Nil".
$ws.Range("I27").Value = "Rejected"
$ws.Range("J27").Value = "Nil"

# Row 28 (TestScenario_27 / TestScenario_27.TestCase_1):
#   Approved/Rejected was "Rejected" with reason "Nil".
#   It is now re-reviewed and marked "Approved" with no reject reason.
$ws.Range("I28").Value = "Approved"
$ws.Range("J28").ClearContents()

# Reflect the resulting view/selection state (row 27 is now the
# selected Approved/Rejected + ReasonToReject pair).
$ws.Range("I27:J27").Select()
